# TrackOrder page implementation / Admincommonbase cleanup regression data
# refresh: bump the 4th data row from the "fake4" fixture to "fake5", and
# nudge the associated phone number + grid selection to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 ("fake4" -> "fake5" test fixture)
$ws.Range("A4").Value = "fake5"
$ws.Range("B4").Value = "fake5@g.com"
$ws.Range("C4").Value = 987654341
$ws.Range("D4").Value = "fake@12345"
$ws.Range("E4").Value = "fake@12345"

# Move the active selection from F10 to F9
$null = $ws.Range("F9").Select()
